# Auto-generated edit script: updates profit/price calculation cells
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets per scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(18, 8).Value = 15284.2
$ws.Cells.Item(18, 9).Value = 10587.375
$ws.Cells.Item(18, 11).Value = 10587.375
$ws.Cells.Item(18, 13).Value = -10303.375

$ws.Cells.Item(19, 8).Value = 1728.6
$ws.Cells.Item(19, 9).Value = 1616
$ws.Cells.Item(19, 10).Value = 1803.6666
$ws.Cells.Item(19, 11).Value = 1616
$ws.Cells.Item(19, 12).Value = 1803.6666
$ws.Cells.Item(19, 13).Value = -1441
$ws.Cells.Item(19, 14).Value = -2153.6666

$ws.Cells.Item(40, 8).Value = 2246
$ws.Cells.Item(40, 9).Value = 2500
$ws.Cells.Item(40, 11).Value = 2500
$ws.Cells.Item(40, 13).Value = -2325

$ws.Cells.Item(132, 8).Value = 1299.6666
$ws.Cells.Item(132, 9).Value = 1299.6666
$ws.Cells.Item(132, 11).Value = 3898.9998
$ws.Cells.Item(132, 13).Value = -1368.9998

$ws.Cells.Item(138, 8).Value = 3228.8096
$ws.Cells.Item(138, 10).Value = 1682.4
$ws.Cells.Item(138, 12).Value = 5047.200000000001
$ws.Cells.Item(138, 14).Value = -15327.2

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 2951.8147
$ws.Cells.Item(32, 9).Value = 2042.0698
$ws.Cells.Item(32, 11).Value = 2042.0698
$ws.Cells.Item(32, 13).Value = -1755.0698

$ws.Cells.Item(61, 8).Value = 3430.5925
$ws.Cells.Item(61, 9).Value = 2714.077
$ws.Cells.Item(61, 10).Value = 4095.9285
$ws.Cells.Item(61, 11).Value = 2714.077
$ws.Cells.Item(61, 12).Value = 4095.9285
$ws.Cells.Item(61, 13).Value = -2502.077
$ws.Cells.Item(61, 14).Value = -4519.9285

$ws.Cells.Item(132, 8).Value = 2982.3513
$ws.Cells.Item(132, 9).Value = 2851.25
$ws.Cells.Item(132, 11).Value = 8553.75
$ws.Cells.Item(132, 13).Value = -6023.75

$ws.Cells.Item(136, 8).Value = 3430.5925
$ws.Cells.Item(136, 9).Value = 2714.077
$ws.Cells.Item(136, 10).Value = 4095.9285
$ws.Cells.Item(136, 11).Value = 8142.231000000001
$ws.Cells.Item(136, 12).Value = 12287.7855
$ws.Cells.Item(136, 13).Value = -5592.231000000001
$ws.Cells.Item(136, 14).Value = -17387.7855

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2163.3333
$ws.Cells.Item(20, 9).Value = 1995.9166
$ws.Cells.Item(20, 10).Value = 2498.1667
$ws.Cells.Item(20, 11).Value = 1995.9166
$ws.Cells.Item(20, 12).Value = 2498.1667
$ws.Cells.Item(20, 13).Value = -1748.9166
$ws.Cells.Item(20, 14).Value = -2992.1667

$ws.Cells.Item(86, 8).Value = 183189.9
$ws.Cells.Item(86, 9).Value = 1565.4445
$ws.Cells.Item(86, 11).Value = 1565.4445
$ws.Cells.Item(86, 13).Value = -442.4445000000001

$ws.Cells.Item(89, 8).Value = 183189.9
$ws.Cells.Item(89, 9).Value = 1565.4445
$ws.Cells.Item(89, 11).Value = 7827.2225
$ws.Cells.Item(89, 13).Value = -2211.2225

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 1755.069
$ws.Cells.Item(31, 9).Value = 951
$ws.Cells.Item(31, 11).Value = 951
$ws.Cells.Item(31, 13).Value = -656

$ws.Cells.Item(34, 8).Value = 1755.069
$ws.Cells.Item(34, 9).Value = 951
$ws.Cells.Item(34, 11).Value = 951
$ws.Cells.Item(34, 13).Value = -749

$ws.Cells.Item(41, 8).Value = 15039.167
$ws.Cells.Item(41, 10).Value = 29000
$ws.Cells.Item(41, 12).Value = 29000
$ws.Cells.Item(41, 14).Value = -29856

$ws.Cells.Item(50, 8).Value = 18000
$ws.Cells.Item(50, 10).Value = 18000
$ws.Cells.Item(50, 12).Value = 18000
$ws.Cells.Item(50, 14).Value = -19250

$ws.Cells.Item(51, 8).Value = 25975
$ws.Cells.Item(51, 9).Value = 5000
$ws.Cells.Item(51, 10).Value = 32966.668
$ws.Cells.Item(51, 11).Value = 5000
$ws.Cells.Item(51, 12).Value = 32966.668
$ws.Cells.Item(51, 13).Value = -4264
$ws.Cells.Item(51, 14).Value = -34438.668

$ws.Cells.Item(58, 8).Value = 2558929.2
$ws.Cells.Item(58, 9).Value = 3953652
$ws.Cells.Item(58, 11).Value = 3953652
$ws.Cells.Item(58, 13).Value = -3953449

$ws.Cells.Item(59, 8).Value = 39266.332
$ws.Cells.Item(59, 10).Value = 39266.332
$ws.Cells.Item(59, 12).Value = 39266.332
$ws.Cells.Item(59, 14).Value = -41556.332

$ws.Cells.Item(60, 8).Value = 4110.8887
$ws.Cells.Item(60, 10).Value = 14999
$ws.Cells.Item(60, 12).Value = 14999
$ws.Cells.Item(60, 14).Value = -16021

$ws.Cells.Item(61, 8).Value = 25975
$ws.Cells.Item(61, 9).Value = 5000
$ws.Cells.Item(61, 10).Value = 32966.668
$ws.Cells.Item(61, 11).Value = 5000
$ws.Cells.Item(61, 12).Value = 32966.668
$ws.Cells.Item(61, 13).Value = -4652
$ws.Cells.Item(61, 14).Value = -33662.668

$ws.Cells.Item(62, 8).Value = 9991.666999999999
$ws.Cells.Item(62, 9).Value = 9991.666999999999
$ws.Cells.Item(62, 11).Value = 9991.666999999999
$ws.Cells.Item(62, 13).Value = -9367.666999999999

$ws.Cells.Item(65, 8).Value = 9991.666999999999
$ws.Cells.Item(65, 9).Value = 9991.666999999999
$ws.Cells.Item(65, 11).Value = 49958.335
$ws.Cells.Item(65, 13).Value = -46838.335

$ws.Cells.Item(107, 8).Value = 477.15
$ws.Cells.Item(107, 9).Value = 402.44446
$ws.Cells.Item(107, 10).Value = 1149.5
$ws.Cells.Item(107, 11).Value = 402.44446
$ws.Cells.Item(107, 12).Value = 1149.5
$ws.Cells.Item(107, 13).Value = 1517.55554
$ws.Cells.Item(107, 14).Value = -4989.5

$ws.Cells.Item(136, 8).Value = 2558929.2
$ws.Cells.Item(136, 9).Value = 3953652
$ws.Cells.Item(136, 11).Value = 11860956
$ws.Cells.Item(136, 13).Value = -11858406

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(113, 8).Value = 125699.25
$ws.Cells.Item(113, 10).Value = 798.7143
$ws.Cells.Item(113, 12).Value = 2396.1429
$ws.Cells.Item(113, 14).Value = -6736.1429

$ws.Cells.Item(131, 8).Value = 9448288
$ws.Cells.Item(131, 9).Value = 250000510
$ws.Cells.Item(131, 10).Value = 14866.941
$ws.Cells.Item(131, 11).Value = 750001530
$ws.Cells.Item(131, 12).Value = 44600.823
$ws.Cells.Item(131, 13).Value = -749996490
$ws.Cells.Item(131, 14).Value = -54680.823

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(10, 8).Value = 3080000
$ws.Cells.Item(10, 9).Value = 3080000
$ws.Cells.Item(10, 10).Value = 0
$ws.Cells.Item(10, 11).Value = 3080000
$ws.Cells.Item(10, 12).Value = 0
$ws.Cells.Item(10, 13).Value = -3079831
$ws.Cells.Item(10, 14).ClearContents()

$ws.Cells.Item(12, 8).Value = 5325626
$ws.Cells.Item(12, 9).Value = 6500000
$ws.Cells.Item(12, 11).Value = 6500000
$ws.Cells.Item(12, 13).Value = -6499860

$ws.Cells.Item(24, 8).Value = 1548307.8
$ws.Cells.Item(24, 10).Value = 11636.363
$ws.Cells.Item(24, 12).Value = 11636.363
$ws.Cells.Item(24, 14).Value = -11982.363

$ws.Cells.Item(29, 8).Value = 62507.75
$ws.Cells.Item(29, 9).Value = 40007
$ws.Cells.Item(29, 11).Value = 40007
$ws.Cells.Item(29, 13).Value = -39717

$ws.Cells.Item(135, 8).Value = 49165.8
$ws.Cells.Item(135, 10).Value = 49165.8
$ws.Cells.Item(135, 12).Value = 49165.8
$ws.Cells.Item(135, 14).Value = -59305.8

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4015.353
$ws.Cells.Item(7, 9).Value = 2515.2856
$ws.Cells.Item(7, 10).Value = 5065.4
$ws.Cells.Item(7, 11).Value = 2515.2856
$ws.Cells.Item(7, 12).Value = 5065.4
$ws.Cells.Item(7, 13).Value = -2403.2856
$ws.Cells.Item(7, 14).Value = -5289.4

$ws.Cells.Item(23, 8).Value = 10000000
$ws.Cells.Item(23, 9).Value = 10000000
$ws.Cells.Item(23, 11).Value = 10000000
$ws.Cells.Item(23, 13).Value = -9999770

$ws.Cells.Item(98, 8).Value = 49425
$ws.Cells.Item(98, 10).Value = 49425
$ws.Cells.Item(98, 12).Value = 49425
$ws.Cells.Item(98, 14).Value = -55415

$ws.Cells.Item(126, 8).Value = 4015.353
$ws.Cells.Item(126, 9).Value = 2515.2856
$ws.Cells.Item(126, 10).Value = 5065.4
$ws.Cells.Item(126, 11).Value = 7545.8568
$ws.Cells.Item(126, 12).Value = 15196.2
$ws.Cells.Item(126, 13).Value = -5075.8568
$ws.Cells.Item(126, 14).Value = -20136.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(21, 8).Value = 15000
$ws.Cells.Item(21, 10).Value = 15000
$ws.Cells.Item(21, 12).Value = 15000
$ws.Cells.Item(21, 14).Value = -15470

$ws.Cells.Item(35, 8).Value = 15000
$ws.Cells.Item(35, 10).Value = 15000
$ws.Cells.Item(35, 12).Value = 15000
$ws.Cells.Item(35, 14).Value = -15580

$ws.Cells.Item(126, 8).Value = 4827.476
$ws.Cells.Item(126, 9).Value = 4113
$ws.Cells.Item(126, 11).Value = 12339
$ws.Cells.Item(126, 13).Value = -9869

$ws.Cells.Item(132, 8).Value = 4420.222
$ws.Cells.Item(132, 9).Value = 3396.6667
$ws.Cells.Item(132, 10).Value = 4932
$ws.Cells.Item(132, 11).Value = 10190.0001
$ws.Cells.Item(132, 12).Value = 14796
$ws.Cells.Item(132, 13).Value = -7660.000100000001
$ws.Cells.Item(132, 14).Value = -19856

$ws.Cells.Item(136, 8).Value = 16837420
$ws.Cells.Item(136, 9).Value = 25254502
$ws.Cells.Item(136, 10).Value = 3259.4546
$ws.Cells.Item(136, 11).Value = 75763506
$ws.Cells.Item(136, 12).Value = 9778.363799999999
$ws.Cells.Item(136, 13).Value = -75760956
$ws.Cells.Item(136, 14).Value = -14878.3638
